$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one weekly record per row (rows 2-55). A new weekly record
# for "Haba" needs to be inserted as the new row 42, pushing the existing
# rows 42-55 down to 43-56 (dimension grows from A1:R55 to A1:R56).
$ws.Rows.Item(42).Insert()

# Columns A, B, C, E, F, G, H, I, N, Q, R are constant for every data row in
# this sheet, so copy them straight from the (still adjacent) row above.
$ws.Cells.Item(42, 1).Value  = $ws.Cells.Item(41, 1).Value2()   # Mercado ID
$ws.Cells.Item(42, 2).Value  = $ws.Cells.Item(41, 2).Value2()   # Mercado
$ws.Cells.Item(42, 3).Value  = $ws.Cells.Item(41, 3).Value2()   # Región
$ws.Cells.Item(42, 5).Value  = $ws.Cells.Item(41, 5).Value2()   # Codreg
$ws.Cells.Item(42, 6).Value  = $ws.Cells.Item(41, 6).Value2()   # Categoría ID
$ws.Cells.Item(42, 7).Value  = $ws.Cells.Item(41, 7).Value2()   # Categoría
$ws.Cells.Item(42, 8).Value  = $ws.Cells.Item(41, 8).Value2()   # Variedad
$ws.Cells.Item(42, 9).Value  = $ws.Cells.Item(41, 9).Value2()   # Calidad
$ws.Cells.Item(42, 14).Value = $ws.Cells.Item(41, 14).Value2()  # Unidad de comercialización
$ws.Cells.Item(42, 17).Value = $ws.Cells.Item(41, 17).Value2()  # Kg o Unidades
$ws.Cells.Item(42, 18).Value = $ws.Cells.Item(41, 18).Value2()  # Clasificación

# New row-specific values for the inserted record.
$ws.Cells.Item(42, 4).Value  = 44511                 # Fecha
$ws.Cells.Item(42, 10).Value = 400                   # Volumen
$ws.Cells.Item(42, 11).Value = 6000                  # Precio mínimo
$ws.Cells.Item(42, 12).Value = 6000                  # Precio máximo
$ws.Cells.Item(42, 13).Value = 6000                  # Precio promedio ponderado
$ws.Cells.Item(42, 15).Value = "Región del Maule"    # Origen
$ws.Cells.Item(42, 16).Value = 240                   # Precio $/Kg
